$wb = $excel.ActiveWorkbook
$wsGeo2 = $wb.Worksheets.Item("geo_id_2_rr_IHD_WHO_2003a")

# Duplicate the geo_id_2 sheet, placing the copy immediately after it.
$wsGeo2.Copy($null, $wsGeo2)
$wsGeo3 = $wb.Worksheets.Item("geo_id_2_rr_IHD_WHO_2003a (2)")
$wsGeo3.Name = "geo_id_3_rr_IHD_WHO_2003a (2)"

# New sheet uses an age offset of 10 (instead of 5 on the original sheet).
$wsGeo3.Range("A12").Value = 10

# Original geo_id_2 sheet: population-fraction baseline offset corrected to 0.
$wsGeo2.Range("A16").Value = 0

foreach ($s in $wb.Worksheets) {
    Write-Host ("Sheet: " + $s.Name + " Index:" + $s.Index)
}
